$d = $word.ActiveDocument

function Find-ReplaceInRange($range, $findText, $replaceText) {
    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. "Exciting opportunity to lead a passionate team" -> "...to be part of..."
# ---------------------------------------------------------------------------
Find-ReplaceInRange $d.Content "Exciting opportunity to lead a passionate team" "Exciting opportunity to be part of a passionate team"

# ---------------------------------------------------------------------------
# 2. "...and are seeking experienced..." -> "...and are now seeking experienced..."
# ---------------------------------------------------------------------------
Find-ReplaceInRange $d.Content "and are seeking experienced" "and are now seeking experienced"

# ---------------------------------------------------------------------------
# 3. First "Essential Skills:" heading -> "Your responsibilities will include:"
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(14).Range.Text = "Your responsibilities will include:"

# ---------------------------------------------------------------------------
# 4. Replace the 4 list items below it with the new responsibilities text
# ---------------------------------------------------------------------------
Find-ReplaceInRange $d.Paragraphs.Item(16).Range "High level understanding of Java, PHP and Python programming languages for back-end development." "Perform and check code reviews."
$d.Paragraphs.Item(17).Range.Text = "Develop U.I and U.X for smartphone applications."
Find-ReplaceInRange $d.Paragraphs.Item(18).Range "Knowledge of API (REST and SOAP)." "Communicate with front-end and back-end team members"
$d.Paragraphs.Item(19).Range.Text = "Receive feedback from clients, testers and stakeholders and perform adjustments."

# ---------------------------------------------------------------------------
# 5. Re-insert a blank paragraph + "Essential Skills:" heading + blank bold
#    paragraph after the new responsibilities list (item 19).
# ---------------------------------------------------------------------------
$p19 = $d.Paragraphs.Item(19)
$p19.Range.InsertParagraphAfter()
$blank1 = $d.Paragraphs.Item(20)
$blank1.Range.ParagraphFormat.Reset()
$blank1.Style = "Normal"
$blank1.Range.ListFormat.RemoveNumbers()

$blank1.Range.InsertParagraphAfter()
$heading = $d.Paragraphs.Item(21)
$heading.Range.ParagraphFormat.Reset()
$heading.Style = "Normal"
$heading.Range.ListFormat.RemoveNumbers()
$heading.Range.Text = "Essential Skills:"
$heading.Range.Bold = 1
$heading.Range.BoldBi = 1

$heading.Range.InsertParagraphAfter()
$blank2 = $d.Paragraphs.Item(22)
$blank2.Range.ParagraphFormat.Reset()
$blank2.Style = "Normal"
$blank2.Range.ListFormat.RemoveNumbers()
$blank2.Range.Bold = 1
$blank2.Range.BoldBi = 1

# ---------------------------------------------------------------------------
# 6. Insert 4 new list paragraphs (copies of the original first 4 items)
#    before "Database handling in SQL." (now item 23)
# ---------------------------------------------------------------------------
$dbItem = $d.Paragraphs.Item(23)
$dbItem.Range.InsertParagraphBefore()
$dbItem.Range.InsertParagraphBefore()
$dbItem.Range.InsertParagraphBefore()
$dbItem.Range.InsertParagraphBefore()

$d.Paragraphs.Item(23).Range.Text = "High level understanding of Java, PHP and Python programming languages for back-end development."
$d.Paragraphs.Item(24).Range.Text = "Professional and proven experience using Android Studio and Apple XCode."
$d.Paragraphs.Item(25).Range.Text = "Knowledge of API (REST and SOAP)."
$d.Paragraphs.Item(26).Range.Text = "Innovative and smart U.I and U.X design."

# ---------------------------------------------------------------------------
# 7. "Experience using C#" -> "Experience using SWIFT"
# ---------------------------------------------------------------------------
Find-ReplaceInRange $d.Paragraphs.Item(31).Range "Experience using C# and developing iOS applications is highly advantageous." "Experience using SWIFT and developing iOS applications is highly advantageous."

# ---------------------------------------------------------------------------
# 8. Before "Application Process:" there used to be 2 blank bold paragraphs;
#    now there should be 4. Insert 2 additional ones.
# ---------------------------------------------------------------------------
$appHeading = $d.Paragraphs.Item(34)
$appHeading.Range.InsertParagraphBefore()
$appHeading.Range.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# 9. "Provide CV/ resume..." -> "Please provide your CV/ resume..."
# ---------------------------------------------------------------------------
Find-ReplaceInRange $d.Content "Provide CV/ resume and cover letter addressing the requirements of the advertised position." "Please provide your CV/ resume and cover letter addressing the requirements of the advertised position."

Write-Host "Paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host "$i : [$($d.Paragraphs.Item($i).Range.Text)]"
}
